# Add two new trailing rows (104 and 105) to each of the 6 worksheets,
# extending the date series by one business day (45967, 45968) and
# appending the corresponding remn_amt values. The second new row's
# amount is 0 for every sheet.

$wb = $excel.ActiveWorkbook

$newAmounts = @{
    1 = @(1303171, 0)
    2 = @(1020734, 0)
    3 = @(507872, 0)
    4 = @(964374, 0)
    5 = @(1724923, 0)
    6 = @(835869, 0)
}

$dateFormat = "YYYY-MM-DD HH:MM:SS"

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $amounts = $newAmounts[$i]

    $ws.Range("A104").Value = 45967
    $ws.Range("A104").NumberFormat = $dateFormat
    $ws.Range("B104").Value = $amounts[0]

    $ws.Range("A105").Value = 45968
    $ws.Range("A105").NumberFormat = $dateFormat
    $ws.Range("B105").Value = $amounts[1]
}
